$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new data rows on Sheet1 (column B, rows 10-15)
$ws1.Range("B10").Value = "hh"
$ws1.Range("B11").Value = "khalif"
$ws1.Range("B12").Value = "KHALIF"
$ws1.Range("B13").Value = "HaYe"
$ws1.Range("B14").Value = "yaah"
$ws1.Range("B15").Value = "Yes"

# Add a new worksheet (Sheet2) right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Select rows 1-9 on the new sheet, then switch back to Sheet1 and
# leave the selection on B15, matching the saved view state.
$ws2.Range("A1:XFD9").Select()
$ws1.Select()
$ws1.Range("B15").Select()
